# Add two new columns, I ("I0") and J ("IF"), to the gausman_kevin sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: new headers take on the same bold/bordered style as the
# existing header cells (e.g. H1).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data rows 2-39: column J mirrors the existing "IP" (column H) value for
# every row, while column I is 1 for every row except row 2, which is 9
# (matching column J there too).
$ipValues = @(1,6,7,6,4,5,6,7,5,6,4,5,5,6,5,5,3,7,7,5,7,8,4,7,6,5,6,9,6,6,8,6,6,7,7,4,3,2)

for ($i = 0; $i -lt $ipValues.Length; $i++) {
    $row = $i + 2
    $jValue = $ipValues[$i]

    if ($row -eq 2) {
        $iValue = 9
        $jValue = 9
    } else {
        $iValue = 1
    }

    $ws.Cells.Item($row, 9).Value = $iValue
    $ws.Cells.Item($row, 10).Value = $jValue
}
